# Update "想去人数" (column F) values across the four sheets of the
# 杭州-漫展信息 workbook to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $Updates) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# 展览 (sheet1)
Set-FValues "展览" @{
    2  = 561
    3  = 260
    4  = 590
    5  = 1403
    6  = 717
    7  = 364
    11 = 6514
    12 = 125
    13 = 33
    15 = 4834
    18 = 5966
    19 = 7662
    22 = 777
    23 = 4097
    24 = 588
    25 = 43
    26 = 75
    27 = 242
    28 = 148
    29 = 1082
    30 = 1519
    32 = 722
    33 = 1722
    34 = 256
    35 = 1974
    37 = 1283
    39 = 714
    40 = 331
    41 = 2495
    42 = 3774
    45 = 460
    49 = 3980
}

# 演出 (sheet2)
Set-FValues "演出" @{
    3  = 1299
    18 = 16
}

# 本地生活 (sheet3)
Set-FValues "本地生活" @{
    2 = 4601
}

# 全部类型 (sheet4)
Set-FValues "全部类型" @{
    2  = 4601
    4  = 561
    5  = 1299
    8  = 260
    9  = 590
    10 = 1403
    11 = 717
    12 = 364
    16 = 6514
    17 = 125
    19 = 4834
    20 = 5966
    21 = 5966
    23 = 777
    24 = 4097
    25 = 588
    26 = 242
    28 = 148
    29 = 1519
    30 = 722
    31 = 1722
    32 = 256
    33 = 1974
    38 = 714
    39 = 331
    41 = 3774
    44 = 460
    49 = 3980
}

$wb.Save()
